$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cells = @("C3","C4","C5","C6","C12")
$styles = @{}
$values = @{}
foreach ($addr in $cells) {
  $styles[$addr] = $ws.Range($addr).Style
  $values[$addr] = $ws.Range($addr).Text
}

# wipe all hyperlinks on the sheet
$ws.Range("A1").Hyperlinks.Delete()

# re-add them; only C5 gets an explicit (different) display text
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:hellothere@tide.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:P@ss1234")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:hellothere@tide.com", "", "", "hellothere@tide.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:P@ss1234")
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:jackjone@tide.com")

foreach ($addr in $cells) {
  $ws.Range($addr).Style = $styles[$addr]
}
foreach ($addr in $cells) {
  if ($addr -eq "C5") {
    $ws.Range($addr).Value = "'hello@tide.com"
  } else {
    $ws.Range($addr).Value = "'" + $values[$addr]
  }
}

Write-Host "count:" $ws.Hyperlinks.Count
